$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# The "Ready for handoff" status text (shared by Overview!E3/F3 and the
# per-locale Status column C3 on both locale sheets) becomes
# "Handback transform failed" everywhere it is used.
[void]$wsOverview.Cells.Replace("Ready for handoff", "Handback transform failed")
[void]$wsZhCn.Cells.Replace("Ready for handoff", "Handback transform failed")
[void]$wsDeDe.Cells.Replace("Ready for handoff", "Handback transform failed")

# Record the Error Detail (column P) for row 3 on each locale sheet, and
# widen that column to fit the longer message.
$wsZhCn.Range("P3").Value = "Handback file name: bwyzvzqh.1ae is different with handoff file name: 6d470eab-8fab-4127-9d7e-fc71f93d4114.2d22762ded34fb1a235f630b17d6508a12e578c3.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

$wsDeDe.Range("P3").Value = "Handback file name: bwyzvzqh.1ae is different with handoff file name: 6d470eab-8fab-4127-9d7e-fc71f93d4114.2d22762ded34fb1a235f630b17d6508a12e578c3.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
